$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.589.41"

$ws.Range("D3").Value = "1.923.72"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").Value = "'326.07"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("D7").Value = "'0.4813"
$ws.Range("E7").Value = "  -0.52%  "

$ws.Range("D8").Value = "'0.4047"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").Value = "'0.08185"
$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("D10").Value = "'1.006"
$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("D11").Value = "'23.74"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "'6.073"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "1.888.74"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("D14").Value = "'7.281"
$ws.Range("E14").Value = "  +1.19%  "

$ws.Range("D15").Value = "'91.40"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").Value = "'0.06870"
$ws.Range("E16").Value = "  +1.79%  "

$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D19").Value = "'17.59"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").Value = "29.588.73"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "'5.655"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D23").Value = "'11.94"
$ws.Range("E23").Value = "  +1.44%  "

$ws.Range("D24").Value = "'2.183"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "2.139.13"
$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").Value = "'156.19"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.01"
$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'6.340"
$ws.Range("E28").Value = "  -3.10%  "

$ws.Range("D29").Value = "'2.087"
$ws.Range("E29").Value = "  -1.62%  "

$ws.Range("D30").Value = "'120.61"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "'1.004"
$ws.Range("E31").Value = "  -2.26%  "

$ws.Range("D32").Value = "'0.09591"
$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").Value = "'5.589"
$ws.Range("E33").Value = "  +1.38%  "

$ws.Range("D34").Value = "'3.559"
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").Value = "'1.390"
$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("D36").Value = "'0.06508"
$ws.Range("E36").Value = "  +6.12%  "

$ws.Range("D37").Value = "'0.02277"
$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("D38").Value = "'1.221"
$ws.Range("E38").Value = "  +3.18%  "

$ws.Range("D39").Value = "'0.5928"
$ws.Range("E39").Value = "  -0.74%  "

$ws.Range("E40").Value = "  -1.26%  "

$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("D42").Value = "'7.845"
$ws.Range("E42").Value = "  -1.88%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1840"
$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.491"
$ws.Range("E44").Value = "  +3.64%  "

$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").Value = "'12.34"
$ws.Range("E46").Value = "  -0.88%  "

$ws.Range("D47").Value = "'0.07513"
$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("D48").Value = "'0.5535"
$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("D49").Value = "'1.947"
$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").Value = "'118.20"
$ws.Range("E50").Value = "  +1.17%  "

$ws.Range("D51").Value = "'2.430"
$ws.Range("E51").Value = "  +0.49%  "
